# hasil ujicoba.xlsx -- "sampek cek detik ter-recover 20 node"
#
# The sheet had a single "Detik ter-recover" (seconds-to-recover) label
# column I used only as a placeholder that actually held the scenario
# sub-label ("Node X"). This change inserts a new column J that takes over
# that sub-label text, freeing column I to hold the real "Detik ter-recover"
# number -- but only for the scenarios that have been measured so far
# (the "10 Node" and "20 Node" sections). The remaining sections (30/40/50
# Node) just get their label text shifted from I to J, with I left blank
# until those seconds-to-recover numbers are gathered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5, 10).Value = 'Node 9'
$ws.Cells.Item(5, 9).Value = 0.099412
$ws.Cells.Item(6, 10).Value = 'Node 9'
$ws.Cells.Item(6, 9).Value = 0.095963
$ws.Cells.Item(8, 10).Value = 'Node 6'
$ws.Cells.Item(8, 9).Value = 0.007739
$ws.Cells.Item(9, 10).Value = 'Node 6'
$ws.Cells.Item(9, 9).Value = 0.008899
$ws.Cells.Item(11, 10).Value = 'Node 4'
$ws.Cells.Item(11, 9).Value = 0.007025
$ws.Cells.Item(12, 10).Value = 'Node 2'
$ws.Cells.Item(12, 9).Value = 0.007022
$ws.Cells.Item(18, 10).Value = 'Node 1'
$ws.Cells.Item(18, 9).Value = 0.035871
$ws.Cells.Item(19, 10).Value = 'Node 7'
$ws.Cells.Item(19, 9).Value = 0.052088
$ws.Cells.Item(20, 10).Value = 'Node 1 + 0'
$ws.Cells.Item(20, 9).Value = 0.030953
$ws.Cells.Item(21, 10).Value = 'Node 0 + 2'
$ws.Cells.Item(21, 9).Value = 0.037118
$ws.Cells.Item(22, 10).Value = 'Node 2 + 0'
$ws.Cells.Item(22, 9).Value = 0.036832
$ws.Cells.Item(24, 10).Value = 'Node 2'
$ws.Cells.Item(24, 9).Value = 0.058625
$ws.Cells.Item(25, 10).Value = 'Node 0'
$ws.Cells.Item(25, 9).Value = 0.040516
$ws.Cells.Item(26, 10).Value = 'Node 2 + 1'
$ws.Cells.Item(26, 9).Value = 0.036499
$ws.Cells.Item(27, 10).Value = 'Node 0 + 3'
$ws.Cells.Item(27, 9).Value = 0.032611
$ws.Cells.Item(28, 10).Value = 'Node 2 + 5'
$ws.Cells.Item(28, 9).Value = 0.041835
$ws.Cells.Item(30, 10).Value = 'Node 8'
$ws.Cells.Item(30, 9).Value = 0.043467
$ws.Cells.Item(31, 10).Value = 'Node 8'
$ws.Cells.Item(31, 9).Value = 0.044429
$ws.Cells.Item(32, 10).Value = 'Node 8 + 3'
$ws.Cells.Item(32, 9).Value = 0.031393
$ws.Cells.Item(33, 10).Value = 'Node 0 + 10'
$ws.Cells.Item(33, 9).Value = 0.034425
$ws.Cells.Item(34, 10).Value = 'Node 2 + 1'
$ws.Cells.Item(34, 9).Value = 0.03065
$ws.Cells.Item(40, 10).Value = 'Node 22'
$ws.Cells.Item(40, 9).ClearContents()
$ws.Cells.Item(41, 10).Value = 'Node 24'
$ws.Cells.Item(41, 9).ClearContents()
$ws.Cells.Item(42, 10).Value = 'Node 26 + 6'
$ws.Cells.Item(42, 9).ClearContents()
$ws.Cells.Item(43, 10).Value = 'Node 13 + 20'
$ws.Cells.Item(43, 9).ClearContents()
$ws.Cells.Item(44, 10).Value = 'Node 24 + 13'
$ws.Cells.Item(44, 9).ClearContents()
$ws.Cells.Item(46, 10).Value = 'Node 16'
$ws.Cells.Item(46, 9).ClearContents()
$ws.Cells.Item(47, 10).Value = 'Node 13'
$ws.Cells.Item(47, 9).ClearContents()
$ws.Cells.Item(48, 10).Value = 'Node 16 + 4'
$ws.Cells.Item(48, 9).ClearContents()
$ws.Cells.Item(49, 10).Value = 'Node 13 + 5'
$ws.Cells.Item(49, 9).ClearContents()
$ws.Cells.Item(50, 10).Value = 'Node 24 + 0'
$ws.Cells.Item(50, 9).ClearContents()
$ws.Cells.Item(52, 10).Value = 'Node 26'
$ws.Cells.Item(52, 9).ClearContents()
$ws.Cells.Item(53, 10).Value = 'Node 16'
$ws.Cells.Item(53, 9).ClearContents()
$ws.Cells.Item(54, 10).Value = 'Node 26 + 20'
$ws.Cells.Item(54, 9).ClearContents()
$ws.Cells.Item(55, 10).Value = 'Node 13 + 15'
$ws.Cells.Item(55, 9).ClearContents()
$ws.Cells.Item(56, 10).Value = 'Node 24 + 12'
$ws.Cells.Item(56, 9).ClearContents()
$ws.Cells.Item(62, 10).Value = 'Node 2'
$ws.Cells.Item(62, 9).ClearContents()
$ws.Cells.Item(63, 10).Value = 'Node 32'
$ws.Cells.Item(63, 9).ClearContents()
$ws.Cells.Item(64, 10).Value = 'Node 2 + 6'
$ws.Cells.Item(64, 9).ClearContents()
$ws.Cells.Item(65, 10).Value = 'Node 32 + 6'
$ws.Cells.Item(65, 9).ClearContents()
$ws.Cells.Item(66, 10).Value = 'Node 2 + 3'
$ws.Cells.Item(66, 9).ClearContents()
$ws.Cells.Item(67, 10).Value = 'Node 2 + 6 + 1'
$ws.Cells.Item(67, 9).ClearContents()
$ws.Cells.Item(68, 10).Value = 'Node 32 + 6 + 7'
$ws.Cells.Item(68, 9).ClearContents()
$ws.Cells.Item(69, 10).Value = 'Node 2 + 6 + 1 + 7'
$ws.Cells.Item(69, 9).ClearContents()
$ws.Cells.Item(70, 10).Value = 'Node 32 + 6 + 7 + 15'
$ws.Cells.Item(70, 9).ClearContents()
$ws.Cells.Item(71, 10).Value = 'Node 2 + 1 + 3 + 15'
$ws.Cells.Item(71, 9).ClearContents()
$ws.Cells.Item(73, 10).Value = 'Node 33'
$ws.Cells.Item(73, 9).ClearContents()
$ws.Cells.Item(74, 10).Value = 'Node 1'
$ws.Cells.Item(74, 9).ClearContents()
$ws.Cells.Item(75, 10).Value = 'Node 33 + 2'
$ws.Cells.Item(75, 9).ClearContents()
$ws.Cells.Item(76, 10).Value = 'Node 1 + 33'
$ws.Cells.Item(76, 9).ClearContents()
$ws.Cells.Item(77, 10).Value = 'Node 33 + 2'
$ws.Cells.Item(77, 9).ClearContents()
$ws.Cells.Item(78, 10).Value = 'Node 33 + 2 + 3'
$ws.Cells.Item(78, 9).ClearContents()
$ws.Cells.Item(79, 10).Value = 'Node 1 + 33 + 3'
$ws.Cells.Item(79, 9).ClearContents()
$ws.Cells.Item(80, 10).Value = 'Node 33 + 2 + 3 + 6'
$ws.Cells.Item(80, 9).ClearContents()
$ws.Cells.Item(81, 10).Value = 'Node 1 + 33 + 3 + 2'
$ws.Cells.Item(81, 9).ClearContents()
$ws.Cells.Item(84, 10).Value = 'Node 32'
$ws.Cells.Item(84, 9).ClearContents()
$ws.Cells.Item(85, 10).Value = 'Node 13'
$ws.Cells.Item(85, 9).ClearContents()
$ws.Cells.Item(86, 10).Value = 'Node 32 + 13'
$ws.Cells.Item(86, 9).ClearContents()
$ws.Cells.Item(87, 10).Value = 'Node 13 + 26'
$ws.Cells.Item(87, 9).ClearContents()
$ws.Cells.Item(88, 10).Value = 'Node 32 + 1'
$ws.Cells.Item(88, 9).ClearContents()
$ws.Cells.Item(89, 10).Value = 'Node 32 + 13 + 3'
$ws.Cells.Item(89, 9).ClearContents()
$ws.Cells.Item(90, 10).Value = 'Node 13 + 26 + 25'
$ws.Cells.Item(90, 9).ClearContents()
$ws.Cells.Item(91, 10).Value = 'Node 32 + 13 + 3 + 9'
$ws.Cells.Item(91, 9).ClearContents()
$ws.Cells.Item(92, 10).Value = 'Node 13 + 26 + 25 + 36'
$ws.Cells.Item(92, 9).ClearContents()
$ws.Cells.Item(93, 10).Value = 'Node 32 + 13 + 1 + 3'
$ws.Cells.Item(93, 9).ClearContents()
$ws.Cells.Item(99, 10).Value = 'Node 13'
$ws.Cells.Item(99, 9).ClearContents()
$ws.Cells.Item(100, 10).Value = 'Node 13'
$ws.Cells.Item(100, 9).ClearContents()
$ws.Cells.Item(101, 10).Value = 'Node 13 + 4'
$ws.Cells.Item(101, 9).ClearContents()
$ws.Cells.Item(102, 10).Value = 'Node 13 + 14'
$ws.Cells.Item(102, 9).ClearContents()
$ws.Cells.Item(103, 10).Value = 'Node 13 + 14'
$ws.Cells.Item(103, 9).ClearContents()
$ws.Cells.Item(104, 10).Value = 'Node 13 + 4 + 5'
$ws.Cells.Item(104, 9).ClearContents()
$ws.Cells.Item(105, 10).Value = 'Node 13 + 14 + 0'
$ws.Cells.Item(105, 9).ClearContents()
$ws.Cells.Item(106, 10).Value = 'Node 13 + 4 + 5 + 6'
$ws.Cells.Item(106, 9).ClearContents()
$ws.Cells.Item(107, 10).Value = 'Node 13 + 14 + 0 + 5'
$ws.Cells.Item(107, 9).ClearContents()
$ws.Cells.Item(108, 10).Value = 'Node 13 + 4 + 14 + 0'
$ws.Cells.Item(108, 9).ClearContents()
$ws.Cells.Item(110, 10).Value = 'Node 4'
$ws.Cells.Item(110, 9).ClearContents()
$ws.Cells.Item(111, 10).Value = 'Node 3'
$ws.Cells.Item(111, 9).ClearContents()
$ws.Cells.Item(112, 10).Value = 'Node 4 + 27'
$ws.Cells.Item(112, 9).ClearContents()
$ws.Cells.Item(113, 10).Value = 'Node 3 + 26'
$ws.Cells.Item(113, 9).ClearContents()
$ws.Cells.Item(114, 10).Value = 'Node 0 + 3'
$ws.Cells.Item(114, 9).ClearContents()
$ws.Cells.Item(115, 10).Value = 'Node 4 + 27 + 36'
$ws.Cells.Item(115, 9).ClearContents()
$ws.Cells.Item(116, 10).Value = 'Node 3 + 26 + 0'
$ws.Cells.Item(116, 9).ClearContents()
$ws.Cells.Item(117, 10).Value = 'Node 4 + 27 + 36 + 44'
$ws.Cells.Item(117, 9).ClearContents()
$ws.Cells.Item(118, 10).Value = 'Node 3 + 26 + 0 + 32'
$ws.Cells.Item(118, 9).ClearContents()
$ws.Cells.Item(119, 10).Value = 'Node 0 + 10 + 3 + 32'
$ws.Cells.Item(119, 9).ClearContents()
$ws.Cells.Item(121, 10).Value = 'Node 1'
$ws.Cells.Item(121, 9).ClearContents()
$ws.Cells.Item(122, 10).Value = 'Node 22'
$ws.Cells.Item(122, 9).ClearContents()
$ws.Cells.Item(123, 10).Value = 'Node 1 + 7'
$ws.Cells.Item(123, 9).ClearContents()
$ws.Cells.Item(124, 10).Value = 'Node 22 + 7'
$ws.Cells.Item(124, 9).ClearContents()
$ws.Cells.Item(125, 10).Value = 'Node 0 + 22'
$ws.Cells.Item(125, 9).ClearContents()
$ws.Cells.Item(126, 10).Value = 'Node 1 + 7 + 10'
$ws.Cells.Item(126, 9).ClearContents()
$ws.Cells.Item(127, 10).Value = 'Node 22 + 7 + 0'
$ws.Cells.Item(127, 9).ClearContents()
$ws.Cells.Item(128, 10).Value = 'Node 1 + 7 + 10 + 0'
$ws.Cells.Item(128, 9).ClearContents()
$ws.Cells.Item(129, 10).Value = 'Node 22 + 7 + 0 + 3'
$ws.Cells.Item(129, 9).ClearContents()
$ws.Cells.Item(130, 10).Value = 'Node 0 + 10 + 22 + 7'
$ws.Cells.Item(130, 9).ClearContents()

# Restore the view state (scrolled/selected cell moved while editing).
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("I34").Select()
